$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 721.75
$ws.Range("I19").Value = 835.3333
$ws.Range("J19").Value = 608.1667
$ws.Range("K19").Value = 835.3333
$ws.Range("L19").Value = 608.1667
$ws.Range("M19").Value = -660.3333
$ws.Range("N19").Value = -958.1667

$ws.Range("H46").Value = 3017
$ws.Range("I46").Value = 3017
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 9051
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -8932
$ws.Range("N46").ClearContents()

$ws.Range("H60").Value = 3017
$ws.Range("I60").Value = 3017
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 9051
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -8567
$ws.Range("N60").ClearContents()

$ws.Range("H101").Value = 2268.5557
$ws.Range("I101").Value = 1278
$ws.Range("J101").Value = 4249.6665
$ws.Range("K101").Value = 3834
$ws.Range("L101").Value = 12748.9995
$ws.Range("M101").Value = -2212
$ws.Range("N101").Value = -15992.9995

$ws.Range("H106").Value = 3143.2856
$ws.Range("I106").Value = 3143.2856
$ws.Range("K106").Value = 3143.2856
$ws.Range("M106").Value = -2512.2856

$ws.Range("H113").Value = 7693
$ws.Range("I113").Value = 6850.5713
$ws.Range("J113").Value = 8146.615
$ws.Range("K113").Value = 6850.5713
$ws.Range("L113").Value = 8146.615
$ws.Range("M113").Value = -3596.5713
$ws.Range("N113").Value = -14654.615

$ws.Range("H131").Value = 6665.706
$ws.Range("I131").Value = 4415.5
$ws.Range("J131").Value = 17166.666
$ws.Range("K131").Value = 13246.5
$ws.Range("L131").Value = 51499.99800000001
$ws.Range("M131").Value = -8206.5
$ws.Range("N131").Value = -61579.99800000001

$ws.Range("H137").Value = 2513.4082
$ws.Range("I137").Value = 1772.5333
$ws.Range("J137").Value = 3683.2104
$ws.Range("K137").Value = 5317.5999
$ws.Range("L137").Value = 11049.6312
$ws.Range("M137").Value = -2767.5999
$ws.Range("N137").Value = -16149.6312

$ws.Range("H138").Value = 2719.9678
$ws.Range("I138").Value = 887.1579
$ws.Range("J138").Value = 3190.554
$ws.Range("K138").Value = 2661.4737
$ws.Range("L138").Value = 9571.662
$ws.Range("M138").Value = 2478.5263
$ws.Range("N138").Value = -19851.662

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 12100
$ws.Range("I45").Value = 8750
$ws.Range("J45").Value = 14333.333
$ws.Range("K45").Value = 8750
$ws.Range("L45").Value = 14333.333
$ws.Range("M45").Value = -8373
$ws.Range("N45").Value = -15087.333

$ws.Range("H122").Value = 2466.15
$ws.Range("I122").Value = 1675.6285
$ws.Range("K122").Value = 5026.8855
$ws.Range("M122").Value = -2576.8855

$ws.Range("H132").Value = 5396.1025
$ws.Range("I132").Value = 5726.636
$ws.Range("J132").Value = 4968.353
$ws.Range("K132").Value = 17179.908
$ws.Range("L132").Value = 14905.059
$ws.Range("M132").Value = -14649.908
$ws.Range("N132").Value = -19965.059

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value = 65000
$ws.Range("J53").Value = 65000
$ws.Range("L53").Value = 65000
$ws.Range("N53").Value = -66148

$ws.Range("H134").Value = 3891.074
$ws.Range("I134").Value = 2705.7368
$ws.Range("K134").Value = 8117.2104
$ws.Range("M134").Value = -5582.2104

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4331.532
$ws.Range("J31").Value = 4335.271
$ws.Range("L31").Value = 4335.271
$ws.Range("N31").Value = -4925.271

$ws.Range("H34").Value = 4331.532
$ws.Range("J34").Value = 4335.271
$ws.Range("L34").Value = 4335.271
$ws.Range("N34").Value = -4739.271

$ws.Range("H122").Value = 258002.28
$ws.Range("I122").Value = 445034.12
$ws.Range("K122").Value = 1335102.36
$ws.Range("M122").Value = -1332652.36

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1860.5
$ws.Range("I114").Value = 900
$ws.Range("J114").Value = 2052.6
$ws.Range("K114").Value = 2700
$ws.Range("L114").Value = 6157.799999999999
$ws.Range("M114").Value = 554
$ws.Range("N114").Value = -12665.8

$ws.Range("H117").Value = 2487.05
$ws.Range("I117").Value = 1562.6666
$ws.Range("J117").Value = 2883.2144
$ws.Range("K117").Value = 4687.9998
$ws.Range("L117").Value = 8649.643199999999
$ws.Range("M117").Value = -1245.9998
$ws.Range("N117").Value = -15533.6432

$ws.Range("H131").Value = 3578.9666
$ws.Range("I131").Value = 5064.143
$ws.Range("K131").Value = 15192.429
$ws.Range("M131").Value = -10152.429

$ws.Range("H134").Value = 8554.666999999999
$ws.Range("I134").Value = 2617.1428
$ws.Range("K134").Value = 7851.428400000001
$ws.Range("M134").Value = -2781.428400000001

$ws.Range("H139").Value = 6099.5366
$ws.Range("I139").Value = 2951.6316
$ws.Range("J139").Value = 8818.182000000001
$ws.Range("K139").Value = 8854.8948
$ws.Range("L139").Value = 26454.546
$ws.Range("M139").Value = -3714.8948
$ws.Range("N139").Value = -36734.546

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 25472.5
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 25472.5
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 25472.5
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -26664.5

$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

$ws.Range("H80").Value = 25730262
$ws.Range("I80").Value = 142233.62
$ws.Range("J80").Value = 66671108
$ws.Range("K80").Value = 142233.62
$ws.Range("L80").Value = 66671108
$ws.Range("M80").Value = -141235.62
$ws.Range("N80").Value = -66673104

$ws.Range("H83").Value = 25730262
$ws.Range("I83").Value = 142233.62
$ws.Range("J83").Value = 66671108
$ws.Range("K83").Value = 711168.1
$ws.Range("L83").Value = 333355540
$ws.Range("M83").Value = -706176.1
$ws.Range("N83").Value = -333365524

$ws.Range("H102").Value = 2160.721
$ws.Range("I102").Value = 1680.9512
$ws.Range("K102").Value = 1680.9512
$ws.Range("M102").Value = -58.95119999999997

$ws.Range("H113").Value = 13650
$ws.Range("I113").Value = 10475
$ws.Range("J113").Value = 20000
$ws.Range("K113").Value = 10475
$ws.Range("L113").Value = 20000
$ws.Range("M113").Value = -8305
$ws.Range("N113").Value = -24340

$ws.Range("H122").Value = 4861.2144
$ws.Range("I122").Value = 2686.7144
$ws.Range("J122").Value = 7035.7144
$ws.Range("K122").Value = 8060.1432
$ws.Range("L122").Value = 21107.1432
$ws.Range("M122").Value = -5610.1432
$ws.Range("N122").Value = -26007.1432

$ws.Range("H126").Value = 7402.8
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 7402.8
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 22208.4
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -27148.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 955.44446
$ws.Range("J46").Value = 700
$ws.Range("L46").Value = 700
$ws.Range("N46").Value = -1076

$ws.Range("H104").Value = 35600
$ws.Range("J104").Value = 35600
$ws.Range("L104").Value = 35600
$ws.Range("N104").Value = -42588

$ws.Range("H132").Value = 4753.484
$ws.Range("I132").Value = 3979.913
$ws.Range("K132").Value = 11939.739
$ws.Range("M132").Value = -9409.739

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 32270.637
$ws.Range("I49").Value = 25000
$ws.Range("J49").Value = 32997.7
$ws.Range("K49").Value = 25000
$ws.Range("L49").Value = 32997.7
$ws.Range("M49").Value = -24770
$ws.Range("N49").Value = -33457.7

$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

$ws.Range("H122").Value = 3076.6287
$ws.Range("I122").Value = 2444.2068
$ws.Range("K122").Value = 7332.6204
$ws.Range("M122").Value = -4882.6204

$ws.Range("H132").Value = 3177.8096
$ws.Range("I132").Value = 1902.3077
$ws.Range("K132").Value = 5706.9231
$ws.Range("M132").Value = -3176.9231
